# Rename TestJframe to GUIPrototype - Changed Backlog
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Product Backlog")
$ws2 = $wb.Worksheets.Item("Sprint Backlog")

# --- Product Backlog sheet ---
# Add component "GUI" to the new backlog item in row 2
$ws1.Range("B2").Value = "GUI"

# --- Sprint Backlog sheet ---
# Update status / description / effort values for the GUI prototype story
$ws2.Range("K3").Value = "open"
$ws2.Range("D2").Value = "Design GUI Prototype"
$ws2.Range("K2").Value = "in progress"
$ws2.Range("K4").Value = "open"
$ws2.Range("K5").Value = "open"
$ws2.Range("K6").Value = "open"

# Reset effort figures to 0 for all sprint backlog rows
$ws2.Range("H2").Value = 0
$ws2.Range("I2").Value = 0

$ws2.Range("H3").Value = 0
$ws2.Range("I3").Value = 0
$ws2.Range("J3").Value = 0

$ws2.Range("H4").Value = 0
$ws2.Range("I4").Value = 0
$ws2.Range("J4").Value = 0

$ws2.Range("H5").Value = 0
$ws2.Range("I5").Value = 0
$ws2.Range("J5").Value = 0

$ws2.Range("H6").Value = 0
$ws2.Range("I6").Value = 0
$ws2.Range("J6").Value = 0

# Grow row 2 to fit the wrapped description text
$ws2.Rows.Item(2).RowHeight = 30

# Restore active sheet to Sprint Backlog and set the selections shown in the workbook
$ws2.Select()
$ws1.Range("B3").Select()
$ws2.Select()
$ws2.Range("D2").Select()
